$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 616-636: append latest NAV history data.
# Column A holds date strings as text (matching the existing column A cells),
# so force text format on that range before assigning values (otherwise Excel
# auto-converts "yyyy-mm-dd" looking strings into date serial numbers).
$ws.Range("A616:A636").NumberFormat = "@"

$rows = @(
    @{r=616; A="2024-08-28"; C=1713.5; D=611.2000122070312; E=1138.300048828125; F=180.9600067138672; G=641.5499877929688; H=22112.82034301758; I=0; J=218.9484263325294},
    @{r=617; A="2024-08-29"; C=1755.650024414062; D=603.6199951171875; E=1132.050048828125; F=179.9400024414062; G=644.2999877929688; H=22256.81030273438; I=0.006511605371146772; J=220.3741320814404},
    @{r=618; A="2024-08-30"; C=1783.050048828125; D=600.3599853515625; E=1127.900024414062; F=178.6199951171875; G=632.0499877929688; H=22286.13012695312; I=0.001317341695415713; J=220.6644401142224},
    @{r=619; A="2024-09-02"; C=1840.550048828125; D=608.5800170898438; E=1111.550048828125; F=177.5399932861328; G=670.2000122070312; H=22669.95040893555; I=0.01722238359894636; J=224.4648077485163},
    @{r=620; A="2024-09-03"; C=1865.599975585938; D=599.9400024414062; E=1114; F=178.4600067138672; G=659.0999755859375; H=22746.27993774414; I=0.003366991432787071; J=225.2205788331677},
    @{r=621; A="2024-09-04"; C=1871.900024414062; D=609; E=1127.900024414062; F=176.0200042724609; G=650.8499755859375; H=22782.94021606445; I=0.001611704349926693; J=225.5835678197662},
    @{r=622; A="2024-09-05"; C=1864.949951171875; D=602.1799926757812; E=1115.150024414062; F=173.4799957275391; G=643.8499755859375; H=22586.02963256836; I=-0.008642896027846764; J=223.6338724975092},
    @{r=623; A="2024-09-06"; C=1857.150024414062; D=597.2999877929688; E=1100; F=169.8500061035156; G=634.7000122070312; H=22350.45025634766; I=-0.01043031378480992; J=221.301301034448},
    @{r=624; A="2024-09-09"; C=1860.449951171875; D=610.3400268554688; E=1104.150024414062; F=168.3300018310547; G=635.2000122070312; H=22404.27005004883; I=0.002407995950143633; J=221.8341936711005},
    @{r=625; A="2024-09-10"; C=1824.5; D=608; E=1113.199951171875; F=169.75; G=637.0499877929688; H=22289.99975585938; I=-0.005100380147810443; J=220.7027549535948},
    @{r=626; A="2024-09-11"; C=1833.150024414062; D=627.6599731445312; E=1112.599975585938; F=165.8800048828125; G=627.2000122070312; H=22284.95007324219; I=-0.0002265447587481507; J=220.6527559012188},
    @{r=627; A="2024-09-12"; C=1854.849975585938; D=645.5999755859375; E=1120.099975585938; F=167.0200042724609; G=651.0999755859375; H=22615.88967895508; I=0.01485036334500268; J=223.9295294994281},
    @{r=628; A="2024-09-13"; C=1894.449951171875; D=646.6500244140625; E=1118.550048828125; F=167.25; G=633.4500122070312; H=22746.35009765625; I=0.005768529142701387; J=225.221273516257},
    @{r=629; A="2024-09-16"; C=1857.599975585938; D=621.0499877929688; E=1115.849975585938; F=163.9600067138672; G=665.9500122070312; H=22506.51992797852; I=-0.01054367705799297; J=222.8466131417117},
    @{r=630; A="2024-09-17"; C=1848.699951171875; D=649.6500244140625; E=1110.949951171875; F=160.6000061035156; G=666.3499755859375; H=22484.49969482422; I=-0.0009783935155129372; J=222.6285814604598},
    @{r=631; A="2024-09-18"; C=1888.199951171875; D=646.7000122070312; E=1079.949951171875; F=158.5599975585938; G=651.7000122070312; H=22442.71960449219; I=-0.00185817300358472; J=222.2148990405636},
    @{r=632; A="2024-09-19"; C=1890.400024414062; D=652.1500244140625; E=1054.449951171875; F=155.25; G=649.5999755859375; H=22292.29992675781; I=-0.006702381903139165; J=220.7255299226262},
    @{r=633; A="2024-09-20"; C=1916.800048828125; D=654.4500122070312; E=1054.599975585938; F=161.4299926757812; G=665.1500244140625; H=22632.26013183594; I=0.01525011803156592; J=224.0916203065262},
    @{r=634; A="2024-09-23"; C=1919.949951171875; D=654.0999755859375; E=1055.25; F=159.5599975585938; G=672; H=22635.46960449219; I=0.0001418096397599883; J=224.1233986584751},
    @{r=635; A="2024-09-24"; C=1904.650024414062; D=646.8499755859375; E=1051.550048828125; F=158.7400054931641; G=675.25; H=22510.13034057617; I=-0.005537294613544976; J=222.8823613703141},
    @{r=636; A="2024-09-25"; C=1928.5; D=633.2999877929688; E=1063.449951171875; F=156.9400024414062; G=667.3499755859375; H=22551.57971191406; I=0.001841365230265907; J=223.292769200981}
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
}
